# Add four new classrooms (H, I, J, K) to the "Classroom Occupancy" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns I, J, K, L (1-indexed columns 9-12)
$headers = @("Classroom H", "Classroom I", "Classroom J", "Classroom K")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 9 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Fill value 1 for rows 2 through 66 in the new columns
for ($row = 2; $row -le 66; $row++) {
    for ($i = 0; $i -lt 4; $i++) {
        $col = 9 + $i
        $ws.Cells.Item($row, $col).Value = 1
    }
}

# Update the selection to match the post-edit state (M1:S1, active cell M1)
$ws.Range("M1:S1").Select()
